$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all Fitness values (column C, rows 2 through 252) to 7569
$ws.Range("C2:C252").Value = 7569
